$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "latest row" number format (currently applied to A8) before
# it gets reassigned, so the newly appended row can reuse it.
$latestRowFormat = $ws.Range("A8").NumberFormat

# Row 8 (A8) currently holds the "latest" style (date-only format).
# Now that a new last row is being appended, row 8 reverts to the regular
# date+time format, matching the format used by rows 2-7.
$ws.Range("A8").NumberFormat = $ws.Range("A7").NumberFormat

# Append the new trading day (2021-11-18) as row 9, reusing the format that
# used to mark the last row (date-only format) for the new last row.
$ws.Range("A9").Value = 44518
$ws.Range("A9").NumberFormat = $latestRowFormat
$ws.Range("B9").Value = -1973.7
